# Add a "comments" column to the "Family" and "clinical values" sheets.
# "Family" sheet also gains a "symptoms" column (mirroring the per-person
# sheets) right before the new "comments" column.

$wb = $excel.ActiveWorkbook

# Pixels-per-character offset baked into the engine's ColumnWidth <-> stored
# width conversion (stored width = ColumnWidth + 5/6, rounded to the nearest
# whole pixel) - used below to land as close as possible on the target
# column widths.
$widthOffset = 5 / 6

# --- Sheet "Family": add "symptoms" (D) and "comments" (E) columns ---
$wsFamily = $wb.Worksheets.Item("Family")

$wsFamily.Range("D1").Value = "symptoms"
$wsFamily.Range("E1").Value = "comments"

$wsFamily.Columns.Item(4).ColumnWidth = 32.6640625 - $widthOffset
$wsFamily.Columns.Item(5).ColumnWidth = 32 - $widthOffset

$wsFamily.Range("E1").Select() | Out-Null

# --- Sheet "clinical values": add "comments" (J) column ---
$wsClinical = $wb.Worksheets.Item("clinical values")

$wsClinical.Range("J1").Value = "comments"

$wsClinical.Columns.Item(10).ColumnWidth = 21.83203125 - $widthOffset

$wsClinical.Range("J1").Select() | Out-Null
